$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "18/03/2023"

$ws.Range("B2").Value = 594.1
$ws.Range("C2").Value = 76

$ws.Range("B3").Value = 60
$ws.Range("C3").Value = 76

$ws.Range("B4").Value = 69

$ws.Range("B5").Value = 471
$ws.Range("C5").Value = 76

$ws.Range("B6").Value = 314
$ws.Range("C6").Value = 76

$ws.Range("B7").Value = 106
$ws.Range("C7").Value = 76

$ws.Range("B8").Value = 92

$ws.Range("B9").Value = 448
$ws.Range("C9").Value = 76

$ws.Range("B10").Value = 57
$ws.Range("C10").Value = 38

$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 38

$ws.Range("B12").Value = 27
$ws.Range("C12").Value = 76
